$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text so values like "1.001" or "7.691" are not
# reinterpreted as numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.604.65'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '1.858.78'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '242.63'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('E6').Value = '  -3.33%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.07611'
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').Value = '0.2998'
$ws.Range('E9').Value = '  +0.95%  '
$ws.Range('D10').Value = '24.68'
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').Value = '0.07728'
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('D12').Value = '1.917.14'
$ws.Range('E12').Value = '  +2.77%  '
$ws.Range('D13').Value = '0.6972'
$ws.Range('E13').Value = '  +0.72%  '
$ws.Range('D14').Value = '5.035'
$ws.Range('E14').Value = '  -0.87%  '
$ws.Range('D15').Value = '83.72'
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('D16').Value = '0.000009931'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('D17').Value = '2.142.26'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').Value = '6.252'
$ws.Range('E18').Value = '  +1.86%  '
$ws.Range('D19').Value = '29.692.77'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D20').Value = '234.97'
$ws.Range('E20').Value = '  -0.70%  '
$ws.Range('D21').Value = '12.62'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = '7.691'
$ws.Range('E23').Value = '  -1.28%  '
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = '155.56'
$ws.Range('E25').Value = '  -1.99%  '
$ws.Range('D26').Value = '0.1402'
$ws.Range('D27').Value = '8.518'
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('D28').Value = '17.79'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').Value = '1.479'
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('D30').Value = '0.05804'
$ws.Range('E30').Value = '  -4.21%  '
$ws.Range('E31').Value = '  -2.21%  '
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('D33').Value = '4.045'
$ws.Range('E33').Value = '  -1.19%  '
$ws.Range('D34').Value = '1.891'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('D36').Value = '0.7233'
$ws.Range('E36').Value = '  -1.68%  '
$ws.Range('D37').Value = '2.589'
$ws.Range('E37').Value = '  -0.86%  '
$ws.Range('D38').Value = '1.260.43'
$ws.Range('E38').Value = '  +4.30%  '
$ws.Range('D39').Value = '2.813'
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('D40').Value = '0.01813'
$ws.Range('E40').Value = '  +1.06%  '
$ws.Range('D41').Value = '0.9065'
$ws.Range('E41').Value = '  -1.00%  '
$ws.Range('D42').Value = '6.166'
$ws.Range('E42').Value = '  -2.53%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '2.059.26'
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '68.23'
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('D46').Value = '101.65'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').Value = '7.347'
$ws.Range('E47').Value = '  -2.06%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.00000000119'
$ws.Range('E48').Value = '  -1.81%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '9.214'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').Value = '0.4060'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '1.716'
$ws.Range('E51').Value = '  +1.92%  '
